$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 59 (shifts old rows 59-71 down to 60-72)
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record
$ws.Cells.Item(59, 1).Value = 3
$ws.Cells.Item(59, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44785
$ws.Cells.Item(59, 5).Value = 5
$ws.Cells.Item(59, 6).Value = 100112035
$ws.Cells.Item(59, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 85
$ws.Cells.Item(59, 11).Value = 14000
$ws.Cells.Item(59, 12).Value = 15000
$ws.Cells.Item(59, 13).Value = 14471
$ws.Cells.Item(59, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(59, 16).Value = 965
$ws.Cells.Item(59, 17).Value = 15
$ws.Cells.Item(59, 18).Value = "Hortaliza"
